$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.214.48"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "2.429.15"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.24%  "
$ws.Range("D9").Value = "2.429.57"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("E11").Value = "  +0.11%  "
$ws.Range("E12").Value = "  -2.45%  "
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.59%  "
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Value = "2.861.68"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "62.152.13"
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").Value = "2.431.29"
$ws.Range("E18").Value = "  +0.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "324.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.85"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.50"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "555.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.33%  "
$ws.Range("D28").Value = "2.544.73"
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "0.0₃0940"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  -3.62%  "
$ws.Range("E33").Value = "  -1.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("E35").Value = "  -2.08%  "
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.79"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.55"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("B40").Value = "EthereumClassic"
$ws.Range("C40").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "150.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.00%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "148.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0533"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("E49").Value = "  +0.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0926"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.43%  "
$ws.Range("E51").Value = "  +0.43%  "
